$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full scraped job-posting table, now including the CompanySize (and DatePosted)
# columns the robot can pull from CareerBuilder listings.
$data = @(
    @("JobTitle","Company","Location","JobType","CompanySize","DatePosted"),
    @("Entry-Level Developer - required to work in office","Vaco Technology","Virginia Beach, VA","Full Time","N/A"),
    @("Entry-level Healthcare Digital Technology Developer","Cognizant Technology","Dallas, TX","Full Time","N/A"),
    @("Entry-level EAS Digital Technology Developer","Cognizant Technology","Dallas, TX","Full Time","N/A"),
    @("Entry Level SQL Developer","HAWAII MAINLAND ADMINISTRATORS L","Tempe, AZ","Full Time","N/A"),
    @("Entry-Level .NET Developer / Application Support","Medline Industries, Inc.","Mundelein, IL","Full Time","N/A"),
    @("Entry Level Unreal Engine C++ Developer","Opex","Moorestown, NJ","Full Time","N/A"),
    @("Java Developer - Recent Grads - Entry Level Positions","Cogent Infotech.","Work From Home","Full Time","N/A"),
    @("Mid to Entry Level Software Developer","Robert Half","Deer Park, TX","Full Time","N/A"),
    @("Developer, Entry Level","Sentinel Technologies Inc.","Downers Grove, IL","Full Time","N/A"),
    @("Jr Web Developer (Entry Level)","Planned Systems International, Inc.","Washington, DC","Contractor","N/A")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Header row is bold.
$ws.Range("A1:F1").Font.Bold = $true

# Resize the columns to fit the newly written content.
$ws.Columns.AutoFit() | Out-Null

# Keep printing sane now that the sheet is much wider.
$ws.PageSetup.Orientation = 1

# Mirror where the user's selection / cursor ended up after the scrape.
$ws.Range("B18").Select()
